# Changes format of occupancy KPIs: Typical arrival/departure times are
# now outputted as text in "hh:mm" format (zero-padded hour, e.g. "08:00")
# instead of plain numbers of hours past midnight.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (arrival) and C (departure) hold weekday times; columns E and F
# hold weekend times. Each numeric "hours past midnight" value gets
# reformatted as a zero-padded "HH:00" text string.
$cols = @("B", "C", "E", "F")

foreach ($row in 2..8) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$row")
        $val = $cell.Value2
        if ($val -ne $null -and $val -ne "") {
            $hour = [int]$val
            $text = "{0:D2}:00" -f $hour
            $cell.Value = $text
        }
    }
}
